$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 (Q18, "For successfully scanning data from input file"):
# grade raised from 8 -> 16 (full marks), and the deduction comment is removed.
$ws.Range("E29").Value = 16
$ws.Range("F29").ClearContents()

# Row 30 (Q19, "For correct and properly aligned output"):
# comment updated to reflect the StringIndexOutOfBoundsException reason.
$ws.Range("F30").Value = "(-4) for no output displayed due to StringIndexOutOfBoundsException"

# Row 37 (Q23, "Compilation errors if any"):
# comment updated to reflect the StringIndexOutOfBoundsException reason.
$ws.Range("F37").Value = "For getting StringIndexOutOfBoundsException while running driver class"

# Leave the final selection on F37, matching the saved cursor position.
$ws.Range("F37").Select()
